# Tripadvisor New Orleans shard 160 update:
#  1. Insert a new "State" column into the hotel_info sheet, right after
#     "Hotel_Name" and before "City", populated with "Louisiana" for the
#     existing hotel row.
#  2. Reorder the sheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Add the "State" column to hotel_info -------------------------------
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City"; insert a fresh column before it so
# "State" lands between "Hotel_Name" (B) and "City" (now shifted to D).
$wsHotel.Columns.Item(3).Insert()

$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder sheet tabs: review_info first, hotel_info second -----------
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))
